# resumes data added in newly created resume folder
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# Row 10: resume file for Thrinath moved to the new "resume" folder and
# converted from a QA pdf naming scheme to a plain docx.
$ws.Range("H10").Value = "D:\\K_Thrinath.docx"

# Row 12: replace the old "Kiran" profile test record with the new
# "Uday" resume entry (new hire data placed into the newly created
# resume folder).
$ws.Range("B12").Value = "yes"
$ws.Range("E12").Value = "udaylewi@gmail.com"
$ws.Range("F12").Value = "VWRheUA4MDc0"
$ws.Range("H12").Value = "D:\\Uday_Resume.docx"

# Update the active view/selection on the DATA sheet to reflect where the
# user was working (matches the saved workbook view state: scrolled so
# column B is left-most visible, with B12 selected).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 2
$excel.Goto($ws.Range("B12"), $true)
$ws.Range("B12").Select()
